$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04558966666666667
$ws.Range("H2").Value = 0.136769
$ws.Range("I2").Value = 0.02375599288687187
$ws.Range("J2").Value = 0.02375599288687187
$ws.Range("M2").Value = 0.1112926666666667
$ws.Range("N2").Value = 0.333878
$ws.Range("O2").Value = 0.01397697460904174
$ws.Range("P2").Value = 0.01397697460904174
$ws.Range("Q2").Value = 0.005073795575777778
$ws.Range("R2").Value = 0.045664160182
$ws.Range("S2").Value = 0.0003320369093923842
$ws.Range("T2").Value = 0.0003320369093923842

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04558966666666667
$ws.Range("H3").Value = 0.136769
$ws.Range("I3").Value = 0.02375599288687187
$ws.Range("J3").Value = 0.02375599288687187
$ws.Range("O3").Value = 0.4165551449121381
$ws.Range("P3").Value = 0.4165551449121381
$ws.Range("Q3").Value = 0.1512141010798889
$ws.Range("R3").Value = 1.360926909719
$ws.Range("S3").Value = 0.009895681059522635
$ws.Range("T3").Value = 0.009895681059522635

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04558966666666667
$ws.Range("H4").Value = 0.136769
$ws.Range("I4").Value = 0.02375599288687187
$ws.Range("J4").Value = 0.02375599288687187
$ws.Range("O4").Value = 0.5694678804788202
$ws.Range("P4").Value = 0.5694678804788201
$ws.Range("Q4").Value = 0.2067231066336667
$ws.Range("R4").Value = 1.860507959703
$ws.Range("S4").Value = 0.01352827491795685
$ws.Range("T4").Value = 0.01352827491795685

$ws.Range("I5").Value = 0.1978186777627204
$ws.Range("J5").Value = 0.1978186777627204
$ws.Range("M5").Value = 0.1112926666666667
$ws.Range("N5").Value = 0.333878
$ws.Range("O5").Value = 0.01397697460904174
$ws.Range("P5").Value = 0.01397697460904174
$ws.Range("Q5").Value = 0.04225003504666666
$ws.Range("R5").Value = 0.38025031542
$ws.Range("S5").Value = 0.002764906636283752
$ws.Range("T5").Value = 0.002764906636283752

$ws.Range("I6").Value = 0.1978186777627204
$ws.Range("J6").Value = 0.1978186777627204
$ws.Range("O6").Value = 0.4165551449121381
$ws.Range("P6").Value = 0.4165551449121381
$ws.Range("S6").Value = 0.08240238798177754
$ws.Range("T6").Value = 0.08240238798177754

$ws.Range("I7").Value = 0.1978186777627204
$ws.Range("J7").Value = 0.1978186777627204
$ws.Range("O7").Value = 0.5694678804788202
$ws.Range("P7").Value = 0.5694678804788201
$ws.Range("S7").Value = 0.1126513831446591
$ws.Range("T7").Value = 0.1126513831446591

$ws.Range("I8").Value = 0.7784253293504076
$ws.Range("J8").Value = 0.7784253293504078
$ws.Range("M8").Value = 0.1112926666666667
$ws.Range("N8").Value = 0.333878
$ws.Range("O8").Value = 0.01397697460904174
$ws.Range("P8").Value = 0.01397697460904174
$ws.Range("Q8").Value = 0.1662557743193333
$ws.Range("R8").Value = 1.496301968874
$ws.Range("S8").Value = 0.0108800310633656
$ws.Range("T8").Value = 0.0108800310633656

$ws.Range("I9").Value = 0.7784253293504076
$ws.Range("J9").Value = 0.7784253293504078
$ws.Range("O9").Value = 0.4165551449121381
$ws.Range("P9").Value = 0.4165551449121381
$ws.Range("S9").Value = 0.3242570758708379
$ws.Range("T9").Value = 0.324257075870838

$ws.Range("I10").Value = 0.7784253293504076
$ws.Range("J10").Value = 0.7784253293504078
$ws.Range("O10").Value = 0.5694678804788202
$ws.Range("P10").Value = 0.5694678804788201
$ws.Range("R10").Value = 60.96425976332101
$ws.Range("S10").Value = 0.4432882224162042
$ws.Range("T10").Value = 0.4432882224162042
